# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    9  = 475
    10 = 6190
    13 = 1000
    14 = 241
    17 = 394
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
